$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: was "Temple Dental" row, now becomes "The Welbeck Clinic" row
$ws.Range("A2").Value = "The Welbeck Clinic"
$ws.Range("B2").Value = "20 Welbeck St, London W1G 8ED, Royaume-Uni"
$ws.Range("C2").Value = "thewelbeckclinic.co.uk"
$ws.Range("D2").Value = "+44 20 7486 8100"

# Row 3: was "The Welbeck Clinic - Cosmetic Dentist" row, now becomes "French Dentist London" row
$ws.Range("A3").Value = "French Dentist London"
$ws.Range("B3").Value = "71 Queen's Gate, South Kensington, London SW7 5JT, Royaume-Uni"
$ws.Range("C3").Value = "drsadone.com"
$ws.Range("D3").Value = "+44 20 7373 6899"
